# "Generate Report for Handback" -- populate the handback columns
# (Latest Target File / Latest Handback File / Latest Handback DateTime)
# on the zh-cn and de-de localization-status sheets, and flip the
# Status column from "Ready for handoff" to "Handed back: in sync with en-US".

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Status column (shared string) flips everywhere it is used ---
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Latest Handback DateTime (column H), was the epoch placeholder ---
$wsZhCn.Range("H2").Value = "2016-03-12 10:34:04"
$wsZhCn.Range("H3").Value = "2016-03-12 10:34:04"

$wsDeDe.Range("H2").Value = "2016-03-12 10:34:10"
$wsDeDe.Range("H3").Value = "2016-03-12 10:34:10"

function Add-HandbackLink($ws, $cellRef, $text, $url) {
    $ws.Hyperlinks.Add($ws.Range($cellRef), $url, "", "", $text)
    $ws.Range($cellRef).Style = "HyperLink"
}

# --- zh-cn: Latest Target File (F) / Latest Handback File (G) ---
Add-HandbackLink $wsZhCn "F2" "733363e9-6fe6-4b0f-bfd0-c61240e0c2e8.md" "https://github.com/OpenLocalizationTest/oltest/blob/617ed65e8e9bd6e7438353231fb7f693644c213c/e2e/733363e9-6fe6-4b0f-bfd0-c61240e0c2e8.md"
Add-HandbackLink $wsZhCn "G2" "733363e9-6fe6-4b0f-bfd0-c61240e0c2e8.d99e715bd5b24972a6c5d86162dff29bb4adf978.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/42f3e44805818a80d408e03569f23e02f5bd2a09/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/733363e9-6fe6-4b0f-bfd0-c61240e0c2e8.d99e715bd5b24972a6c5d86162dff29bb4adf978.zh-cn.xlf"
Add-HandbackLink $wsZhCn "F3" "9e25ea47-4442-4835-bbd2-4ea93cfa2490.md" "https://github.com/OpenLocalizationTest/oltest/blob/617ed65e8e9bd6e7438353231fb7f693644c213c/e2e/9e25ea47-4442-4835-bbd2-4ea93cfa2490.md"
Add-HandbackLink $wsZhCn "G3" "9e25ea47-4442-4835-bbd2-4ea93cfa2490.8fec50f1fa6e86d4ecf2bc29c592b85ec1e67829.zh-cn.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/42f3e44805818a80d408e03569f23e02f5bd2a09/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/9e25ea47-4442-4835-bbd2-4ea93cfa2490.8fec50f1fa6e86d4ecf2bc29c592b85ec1e67829.zh-cn.xlf"

# --- de-de: Latest Target File (F) / Latest Handback File (G) ---
Add-HandbackLink $wsDeDe "F2" "733363e9-6fe6-4b0f-bfd0-c61240e0c2e8.md" "https://github.com/OpenLocalizationTest/oltest/blob/617ed65e8e9bd6e7438353231fb7f693644c213c/e2e/733363e9-6fe6-4b0f-bfd0-c61240e0c2e8.md"
Add-HandbackLink $wsDeDe "G2" "733363e9-6fe6-4b0f-bfd0-c61240e0c2e8.d99e715bd5b24972a6c5d86162dff29bb4adf978.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6ff0cf49720a43cd3c97042d8bad9bcbbba3de8e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/733363e9-6fe6-4b0f-bfd0-c61240e0c2e8.d99e715bd5b24972a6c5d86162dff29bb4adf978.de-de.xlf"
Add-HandbackLink $wsDeDe "F3" "9e25ea47-4442-4835-bbd2-4ea93cfa2490.md" "https://github.com/OpenLocalizationTest/oltest/blob/617ed65e8e9bd6e7438353231fb7f693644c213c/e2e/9e25ea47-4442-4835-bbd2-4ea93cfa2490.md"
Add-HandbackLink $wsDeDe "G3" "9e25ea47-4442-4835-bbd2-4ea93cfa2490.8fec50f1fa6e86d4ecf2bc29c592b85ec1e67829.de-de.xlf" "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6ff0cf49720a43cd3c97042d8bad9bcbbba3de8e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/9e25ea47-4442-4835-bbd2-4ea93cfa2490.8fec50f1fa6e86d4ecf2bc29c592b85ec1e67829.de-de.xlf"
